$p = $ppt.ActivePresentation
$p.Slides.Item(25).Delete()
$p.Slides.Item(24).Delete()
$p.Slides.Item(23).Delete()
